# Add new submissions (rows 10-13) to the tracker table on sheet "Tabelle2" range.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (Name Ramp) values first, in the order they were originally typed ---
$ws.Range("C10").Value = "9_JM_MS"
$ws.Range("C11").Value = "10_JM_MS"
$ws.Range("C12").Value = "11_JM_MS"

# --- Column B (Name) values, typed afterwards in a different order ---
$ws.Range("B11").Value = "221123_xgb_reduced_2"
$ws.Range("B12").Value = "221123_ligthgbm"
$ws.Range("B10").Value = "221122_xgb_reduced"

# --- Last new row's Name Ramp value, typed last ---
$ws.Range("C13").Value = "12_JM_MS"

# --- Column A (Date) values: copy number formatting from an existing date cell
# so the new cells reuse the workbook's existing date style instead of Excel
# minting a brand-new (duplicate) number format. ---
$ws.Range("A4").Copy()
$ws.Range("A10:A13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A10").Value = (Get-Date -Year 2022 -Month 11 -Day 22 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("A11").Value = (Get-Date -Year 2022 -Month 11 -Day 23 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("A12").Value = (Get-Date -Year 2022 -Month 11 -Day 23 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("A13").Value = (Get-Date -Year 2022 -Month 11 -Day 23 -Hour 0 -Minute 0 -Second 0).Date

# --- Column D (Hand in) "TRUE" values: typing TRUE directly would store a
# real boolean, but the tracker stores it as the literal text "TRUE" (reusing
# the existing shared string). Build it as a text formula result in a scratch
# cell, then paste-special the *value* into place so it lands as text. ---
$scratch = $ws.Range("Z1")
$scratch.Formula = '="TRUE"'
foreach ($r in 10..13) {
    $scratch.Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4163)
}
$scratch.ClearContents()
$ws.Application.CutCopyMode = $false

# --- Column E (By) values ---
$ws.Range("E10").Value = "Maria"
$ws.Range("E11").Value = "Joao"
$ws.Range("E12").Value = "Maria"
$ws.Range("E13").Value = "Joao"

# --- Grow the "Tabelle2" structured table so the new rows belong to it,
# matching how Excel auto-expands a table when data is entered right below
# it. ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E13"))

# Move the active selection the way Excel leaves it after entering a new row
# right below the freshly added table data.
[void]$ws.Range("A14").Select()
